# Sprint 6 report-sheet update
# - Rename first sheet from "report-sheet" to "plan"
# - Update burndown chart title to reference Sprint 6
# - Add a dashed line style to the "plan" series of the burndown chart
# - Update saved selections on the plan and burndown sheets
# - Update the zoom level on the burndown sheet

$wb = $excel.ActiveWorkbook

$wsPlan = $wb.Worksheets.Item(1)
$wsExecution = $wb.Worksheets.Item(2)
$wsBurndown = $wb.Worksheets.Item(3)

# Rename the "report-sheet" sheet to "plan"
$wsPlan.Name = "plan"

# Update the burndown chart: title text and the "plan" series dash style
$chart = $wsBurndown.ChartObjects(1).Chart
$chart.ChartTitle.Text = "Sprint 6: Burndown Chart Provisional vs Actual"
$chart.SeriesCollection(1).Border.DashStyle = 4

# Restore the selection on the "plan" sheet (keep it the active tab)
$wsPlan.Activate() | Out-Null
$wsPlan.Range("C15").Select() | Out-Null

# Update the zoom level and selection on the "burndown" sheet without
# leaving it as the active tab
$wsBurndown.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$wsBurndown.Range("C34").Select() | Out-Null

# Restore "plan" as the active/selected sheet
$wsPlan.Activate() | Out-Null
